$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New precision/recall/f1 data in columns I:K for several rows (P11 prompts etc.) ---
$ws.Range("I7").Value = 0.57991552999999996
$ws.Range("J7").Value = 0.30175438999999998
$ws.Range("K7").Value = 0.27339846000000001

$ws.Range("I8").Value = 0.48737355999999998
$ws.Range("J8").Value = 0.31578947000000002
$ws.Range("K8").Value = 0.29975924999999998

$ws.Range("I9").Value = 0.48597420000000002
$ws.Range("J9").Value = 0.28671679
$ws.Range("K9").Value = 0.25912464000000002

$ws.Range("I10").Value = 0.51644239000000003
$ws.Range("J10").Value = 0.31528822000000001
$ws.Range("K10").Value = 0.30435465

$ws.Range("I12").Value = 0.55299197
$ws.Range("J12").Value = 0.36892230576441098
$ws.Range("K12").Value = 0.36742376660801201

$ws.Range("I13").Value = 0.53707185999999996
$ws.Range("J13").Value = 0.34486215999999997
$ws.Range("K13").Value = 0.33676861000000002

$ws.Range("I19").Value = 0.51212380999999996
$ws.Range("J19").Value = 0.27669173000000002
$ws.Range("K19").Value = 0.24118174000000001

$ws.Range("I20").Value = 0.49125333999999998
$ws.Range("J20").Value = 0.30476189999999997
$ws.Range("K20").Value = 0.28231458999999998

# Apply the same numeric style these cells' neighbours use (precision/recall/f1 columns)
$ws.Range("I3:K3").Copy()
$ws.Range("I7:K10").PasteSpecial(-4122)
$ws.Range("I12:K13").PasteSpecial(-4122)
$ws.Range("I19:K20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New annotation cell: red-text label for the control-classes row ---
$ws.Range("L23").Value = "<=> control classes"
$ws.Range("L23").Font.Color = 255

# --- Highlight prompt-label cells with orange fill ---
$ws.Range("G5").Interior.Color = 49407
$ws.Range("G21").Interior.Color = 49407
$ws.Range("G22").Interior.Color = 49407
$ws.Range("G24").Interior.Color = 49407

# --- Restore cursor/selection to match the reviewed region ---
$ws.Activate()
$ws.Range("G24,G21:G22").Select()

Write-Host "edit complete"
